# Apply updated "想去人数" (want-to-go count) figures to the
# "展览" (Exhibition) and "全部类型" (All types) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 12970
$ws1.Range("F6").Value  = 91
$ws1.Range("F7").Value  = 50
$ws1.Range("F10").Value = 12943
$ws1.Range("F12").Value = 43
$ws1.Range("F13").Value = 8709
$ws1.Range("F14").Value = 7719
$ws1.Range("F15").Value = 202
$ws1.Range("F19").Value = 989
$ws1.Range("F20").Value = 10
$ws1.Range("F23").Value = 185
$ws1.Range("F24").Value = 326

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 12970
$ws4.Range("F7").Value  = 91
$ws4.Range("F8").Value  = 50
$ws4.Range("F11").Value = 12944
$ws4.Range("F13").Value = 43
$ws4.Range("F14").Value = 8709
$ws4.Range("F15").Value = 7719
$ws4.Range("F16").Value = 202
$ws4.Range("F20").Value = 989
$ws4.Range("F21").Value = 10
$ws4.Range("F26").Value = 185
$ws4.Range("F27").Value = 326
